# BIS-768: Fixed XLS export test files
# Adds the "Unique" column header (L4) to the sample-type export sheet,
# matching the style of the existing "Multivalued" header (K4), and
# updates the active selection to L5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "Multivalued" header cell (K4) onto
# the new L4 cell so the new header matches the bold header style, then
# set its text.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = "Unique"

# Match the workbook's recorded active cell/selection after the edit.
[void]$ws.Range("L5").Select()
